# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the d56253df-ed40-41db-8663-8ed10e624083 file row on the
# locale sheets (zh-cn row 3, de-de row 3), and roll the newest of those
# handoff timestamps up into the Overview sheet's "Latest HO Xliff
# Generate Date" column for that same file (row 3).

$wb = $excel.ActiveWorkbook

# zh-cn: row 3 is the d56253df-ed40-41db-8663-8ed10e624083... file.
#   H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-12 04:50:41"
$wsZhCn.Range("K3").Value = "2016-08-12 04:50:58"

# de-de: row 3 is the same d56253df-ed40-41db-8663-8ed10e624083... file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-12 04:50:46"
$wsDeDe.Range("K3").Value = "2016-08-12 04:51:09"

# Overview: row 3 corresponds to the same file; its "Latest HO Xliff
# Generate Date" column picks up the newest handoff generation date.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-12 04:50:46"
